$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new data row (row 8) for the "SortHeap" algorithm result ---
# Clone formatting from the last existing data row (row 7) so the new
# row carries the same font / border / alignment as the rest of the table.
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "SortHeap"
$ws.Range("C8").Value = 0.00398993492126464
$ws.Rows.Item(8).RowHeight = 15.75

# --- Column C (the "Speed" values) switches from a 4-decimal to a
#     3-decimal number format ---
$ws.Range("C2:C8").NumberFormat = "0.000"

# The newly-added row's speed cell loses the "shrink to fit" alignment
# once it is sorted into place (matches the target formatting).
$ws.Range("C8").ShrinkToFit = $false

# --- Re-sort the B2:C8 block ascending by the Speed column, same as the
#     author re-running Data > Sort after adding the new algorithm ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C2:C8"))
$ws.Sort.SetRange($ws.Range("B2:C8"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# --- Restore the cursor/selection to where the author left it ---
[void]$ws.Range("C11").Select()
